$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings: volume/issue number + report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/8/2025  Through  9/14/2025"

# --- Weekly crime-stat table updates (rows 14-30, columns C:N) ---
# Row 14
$ws.Range("C14").NumberFormat = "General"
$ws.Range("C14").Value = "0"
$ws.Range("L14").Value = -20
$ws.Range("N14").Value = -55.555555555555

# Row 15
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 22
$ws.Range("J15").Value = 27
$ws.Range("K15").Value = -18.518518518518
$ws.Range("L15").Value = 10
$ws.Range("M15").Value = 29.411764705882
$ws.Range("N15").Value = -12

# Row 16
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 30
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = 36.363636363636
$ws.Range("I16").Value = 207
$ws.Range("J16").Value = 178
$ws.Range("K16").Value = 16.292134831460
$ws.Range("L16").Value = 1.970443349753
$ws.Range("M16").Value = 2.475247524752
$ws.Range("N16").Value = -55.863539445629

# Row 17
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 30
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = 19.444444444444
$ws.Range("I17").Value = 381
$ws.Range("J17").Value = 348
$ws.Range("K17").Value = 9.482758620689
$ws.Range("L17").Value = 31.379310344827
$ws.Range("M17").Value = 111.666666666667
$ws.Range("N17").Value = 73.181818181818

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 165
$ws.Range("J18").Value = 122
$ws.Range("K18").Value = 35.245901639344
$ws.Range("L18").Value = -7.303370786516
$ws.Range("M18").Value = -36.046511627907
$ws.Range("N18").Value = -84.149855907781

# Row 19
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 10.526315789473
$ws.Range("G19").Value = 75
$ws.Range("H19").Value = 1.333333333333
$ws.Range("I19").Value = 565
$ws.Range("J19").Value = 623
$ws.Range("K19").Value = -9.309791332263
$ws.Range("L19").Value = 24.175824175824
$ws.Range("M19").Value = 87.707641196013
$ws.Range("N19").Value = 32.009345794392

# Row 20
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 60
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 34
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 382
$ws.Range("J20").Value = 331
$ws.Range("K20").Value = 15.407854984894
$ws.Range("L20").Value = 1.058201058201
$ws.Range("M20").Value = 137.267080745342
$ws.Range("N20").Value = -70.839694656488

# Row 21
$ws.Range("C21").Value = 57
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = 35.714285714285
$ws.Range("F21").Value = 197
$ws.Range("G21").Value = 179
$ws.Range("H21").Value = 10.055865921787
$ws.Range("I21").Value = 1726
$ws.Range("J21").Value = 1632
$ws.Range("K21").Value = 5.759803921568
$ws.Range("L21").Value = 12.884238064094
$ws.Range("M21").Value = 53.014184397163
$ws.Range("N21").Value = -50.713877784123

# Row 22
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -87.5
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = -33.333333333333
$ws.Range("M22").Value = -7.692307692307

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 8
$ws.Range("H23").Value = -27.272727272727
$ws.Range("I23").Value = 78
$ws.Range("J23").Value = 78
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -14.285714285714
$ws.Range("M23").Value = 65.957446808510

# Row 24
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 105.555555555556
$ws.Range("F24").Value = 138
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = 94.366197183098
$ws.Range("I24").Value = 1248
$ws.Range("J24").Value = 921
$ws.Range("K24").Value = 35.504885993485
$ws.Range("L24").Value = 11.727842435094
$ws.Range("M24").Value = 91.705069124424

# Row 25
$ws.Range("C25").Value = 17
$ws.Range("E25").Value = 142.857142857143
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 190.909090909091
$ws.Range("I25").Value = 501
$ws.Range("J25").Value = 348
$ws.Range("K25").Value = 43.965517241379
$ws.Range("L25").Value = 8.207343412527

# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -16.666666666666
$ws.Range("F26").Value = 38
$ws.Range("H26").Value = -19.148936170212
$ws.Range("I26").Value = 481
$ws.Range("J26").Value = 420
$ws.Range("K26").Value = 14.523809523809
$ws.Range("L26").Value = 28.609625668449
$ws.Range("M26").Value = 10.829493087557

# Row 27
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 26
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = -16.129032258064
$ws.Range("L27").Value = -7.142857142857

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 50
$ws.Range("J28").Value = 44
$ws.Range("K28").Value = 13.636363636363
$ws.Range("L28").Value = 2.040816326530

# Row 29
$ws.Range("C29").NumberFormat = "General"
$ws.Range("C29").Value = "0"
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 17
$ws.Range("K29").Value = -35.294117647058
$ws.Range("N29").Value = -62.068965517241

# Row 30
$ws.Range("C30").NumberFormat = "General"
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -75
$ws.Range("J30").Value = 15
$ws.Range("K30").Value = -60
$ws.Range("N30").Value = -76

